# Re-save style formatting fix:
#  - select the whole sheet (as Excel does before a "fit columns to data" pass)
#  - best-fit the column widths: narrow "label" column A, and the uniform
#    5-character-wide year columns B:S
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select entire sheet (mirrors the saved <selection sqref="A1:XFD1048576"/>)
[void]$ws.Cells.Select()

# Column A ("Institution name" labels) - best-fit narrow width
$ws.Columns.Item(1).ColumnWidth = 3.5

# Columns B:S (the 1996-2013 year/count columns) - best-fit to "5"-ish width
$ws.Range("B1:S1").EntireColumn.ColumnWidth = 4.166666666666667
